$wb = $excel.ActiveWorkbook

# Sheet 1: "Exhibition" (展览) - update column F ("want to go" counts)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 71
$ws1.Range("F5").Value = 255
$ws1.Range("F7").Value = 82
$ws1.Range("F12").Value = 105
$ws1.Range("F13").Value = 2355
$ws1.Range("F14").Value = 58
$ws1.Range("F15").Value = 37
$ws1.Range("F16").Value = 522
$ws1.Range("F17").Value = 539
$ws1.Range("F18").Value = 161
$ws1.Range("F20").Value = 44
$ws1.Range("F21").Value = 48
$ws1.Range("F22").Value = 1834
$ws1.Range("F23").Value = 3961
$ws1.Range("F24").Value = 30
$ws1.Range("F25").Value = 59
$ws1.Range("F26").Value = 1178
$ws1.Range("F27").Value = 226
$ws1.Range("F28").Value = 2083
$ws1.Range("F32").Value = 106
$ws1.Range("F33").Value = 288
$ws1.Range("F34").Value = 417
$ws1.Range("F36").Value = 688
$ws1.Range("F37").Value = 435
$ws1.Range("F38").Value = 411

# Sheet 4: "All Types" (全部类型) - update column F ("want to go" counts)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 71
$ws4.Range("F5").Value = 255
$ws4.Range("F7").Value = 82
$ws4.Range("F12").Value = 105
$ws4.Range("F13").Value = 2355
$ws4.Range("F14").Value = 58
$ws4.Range("F16").Value = 37
$ws4.Range("F17").Value = 522
$ws4.Range("F18").Value = 539
$ws4.Range("F19").Value = 161
$ws4.Range("F21").Value = 44
$ws4.Range("F22").Value = 48
$ws4.Range("F23").Value = 1834
$ws4.Range("F24").Value = 3961
$ws4.Range("F25").Value = 30
$ws4.Range("F26").Value = 59
$ws4.Range("F27").Value = 1178
$ws4.Range("F28").Value = 226
$ws4.Range("F29").Value = 2083
$ws4.Range("F33").Value = 106
$ws4.Range("F34").Value = 288
$ws4.Range("F35").Value = 417
$ws4.Range("F37").Value = 688
$ws4.Range("F38").Value = 435
$ws4.Range("F39").Value = 411
